$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New trailing columns added to the header row (GEL / SMT / TAHUN)
$ws.Range("J1").Value = "GEL"
$ws.Range("K1").Value = "SMT"
$ws.Range("L1").Value = "TAHUN"

# Match the saved view state from the diff: selection moved to M6 and the
# window scrolled so column G is the left-most visible column.
$ws.Range("M6").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
